# B6-PowerPoint.pptx — Fri, May 01, 2020  1:06:27 AM
#
# 1) The three data tables (slides 14, 15, 16) are switched from the
#    deck's default "Table_0" style to the built-in "No Style, Table
#    Grid" style ({B18E7987-DCD1-4129-AD9D-4C8D9EE1CDE6}).
# 2) The presentation's theme colour palette is swapped from the
#    "Integral" (Red Violet) scheme to the standard "Office" scheme
#    (the font scheme and format scheme are identical between the two
#    themes, so only the colour values actually move).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$tableSlides = @(14, 15, 16)
foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{B18E7987-DCD1-4129-AD9D-4C8D9EE1CDE6}")
        }
    }
}

# --- 2. Swap the theme colour scheme (Integral -> Office) -----------
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = 0            # dk1      000000
$cs.Colors(2).RGB  = 16777215     # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388      # dk2      44546A
$cs.Colors(4).RGB  = 15132391     # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939     # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501      # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845     # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407        # accent4  FFC000
$cs.Colors(9).RGB  = 12874308     # accent5  4472C4
$cs.Colors(10).RGB = 4697456      # accent6  70AD47
$cs.Colors(11).RGB = 12673797     # hlink    0563C1
$cs.Colors(12).RGB = 7491477      # folHlink 954F72
